# Update cryptos list with latest prices/volumes scraped on
# Sat Jan 20 08:54:54 UTC 2024 (GitHub Actions run).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells are free-form text (thousands separated by
# '.', trailing zeros, sub-script ₃ digit grouping, etc.) so force the
# cells to Text before writing, then restore the default (Normal) style
# so no new persistent number-format is left behind.
$priceCells = @("D2","D3","D4","D5","D6","D9","D10","D13","D14","D16","D17","D18","D19","D20","D21","D22","D23","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D43","D45","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.525.70"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "2.473.58"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "313.14"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "91.72"
$ws.Range("E6").Value = "  -2.94%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +2.35%  "

$ws.Range("D10").Value = "32.44"
$ws.Range("E10").Value = "  -3.20%  "

$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "2.855.65"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("D14").Value = "16.26"
$ws.Range("E14").Value = "  +8.17%  "

$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("D16").Value = "2.452.12"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("D17").Value = "0.770"
$ws.Range("E17").Value = "  -2.40%  "

$ws.Range("D18").Value = "41.519.81"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").Value = "6.52"
$ws.Range("E19").Value = "  +3.11%  "

$ws.Range("D20").Value = "0.0₃0945"
$ws.Range("E20").Value = "  +2.22%  "

$ws.Range("D21").Value = "71.82"
$ws.Range("E21").Value = "  +4.92%  "

$ws.Range("D22").Value = "11.05"
$ws.Range("E22").Value = "  -1.61%  "

$ws.Range("D23").Value = "235.97"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("D27").Value = "24.84"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").Value = "9.68"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").Value = "35.56"
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("D31").Value = "156.10"
$ws.Range("E31").Value = "  +2.98%  "

$ws.Range("D32").Value = "5.45"
$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").Value = "0.0757"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("D35").Value = "17.30"
$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  -8.49%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -5.79%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +2.71%  "

$ws.Range("E39").Value = "  -3.17%  "

$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("E41").Value = "  -4.98%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "1.958.32"
$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("D45").Value = "18.72"
$ws.Range("E45").Value = "  -3.86%  "

$ws.Range("E46").Value = "  -2.84%  "

$ws.Range("D47").Value = "9.03"
$ws.Range("E47").Value = "  +4.28%  "

$ws.Range("D48").Value = "2.714.40"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "97.70"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("D50").Value = "67.26"
$ws.Range("E50").Value = "  -3.34%  "

$ws.Range("D51").Value = "72.01"
$ws.Range("E51").Value = "  -3.48%  "

# Restore default styling on the price cells we text-formatted above.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
